$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '59.888.33'
$ws.Range("E2").Value = '  -1.57%  '
$ws.Range("D3").Value = '2.303.89'
$ws.Range("E3").Value = '  -2.75%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '541.30'
$ws.Range("E5").Value = '  -1.02%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '128.75'
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("E8").Value = '  -3.21%  '
$ws.Range("D9").Value = '2.301.43'
$ws.Range("E9").Value = '  -2.76%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.101'
$ws.Range("E10").Value = '  -0.93%  '
$ws.Range("E11").Value = '  +0.23%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.149'
$ws.Range("E12").Value = '  -0.81%  '
$ws.Range("E13").Value = '  -1.77%  '
$ws.Range("B14").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C14").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D14").Value = '2.717.36'
$ws.Range("E14").Value = '  -2.68%  '
$ws.Range("B15").Value = 'WrappedBTC'
$ws.Range("C15").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D15").Value = '59.831.60'
$ws.Range("E15").Value = '  -1.57%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '23.12'
$ws.Range("E16").Value = '  -4.12%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000132'
$ws.Range("E17").Value = '  -1.64%  '
$ws.Range("D18").Value = '2.300.54'
$ws.Range("E18").Value = '  -3.51%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.41'
$ws.Range("E19").Value = '  -3.21%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '310.26'
$ws.Range("E20").Value = '  -2.32%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.00'
$ws.Range("E21").Value = '  -4.69%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.46'
$ws.Range("E22").Value = '  -7.69%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.999'
$ws.Range("E23").Value = '  -0.04%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '63.44'
$ws.Range("E24").Value = '  -0.11%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.169'
$ws.Range("E25").Value = '  -2.21%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  +0.04%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.74'
$ws.Range("E27").Value = '  -3.72%  '
$ws.Range("E28").Value = '  -1.14%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '171.68'
$ws.Range("E29").Value = '  +0.06%  '
$ws.Range("E30").Value = '  +3.33%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.70'
$ws.Range("E31").Value = '  -2.99%  '
$ws.Range("D32").Value = '0.0₃0713'
$ws.Range("E32").Value = '  -3.69%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.79'
$ws.Range("E33").Value = '  -1.99%  '
$ws.Range("B34").Value = 'ImmutableX'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.35'
$ws.Range("E34").Value = '  -4.66%  '
$ws.Range("B35").Value = 'PolygonEcosystemToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.378'
$ws.Range("E35").Value = '  -1.42%  '
$ws.Range("E36").Value = '  +0.00%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '17.72'
$ws.Range("E37").Value = '  -2.12%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.00'
$ws.Range("E38").Value = '  +0.09%  '
$ws.Range("E39").Value = '  -3.94%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '310.53'
$ws.Range("E40").Value = '  -2.60%  '
$ws.Range("B41").Value = 'OKB'
$ws.Range("C41").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '37.92'
$ws.Range("E41").Value = '  -1.06%  '
$ws.Range("B42").Value = 'Stacks'
$ws.Range("C42").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.50'
$ws.Range("E42").Value = '  -3.09%  '
$ws.Range("B43").Value = 'Aave'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '136.22'
$ws.Range("E43").Value = '  -5.22%  '
$ws.Range("B44").Value = 'Filecoin'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.40'
$ws.Range("E44").Value = '  -2.12%  '
$ws.Range("B45").Value = 'Stellar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0934'
$ws.Range("E45").Value = '  -2.24%  '
$ws.Range("B46").Value = 'Mantle'
$ws.Range("C46").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.564'
$ws.Range("E46").Value = '  -0.20%  '
$ws.Range("B47").Value = 'InjectiveProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '18.53'
$ws.Range("E47").Value = '  -5.01%  '
$ws.Range("B48").Value = 'Hedera'
$ws.Range("C48").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0488'
$ws.Range("E48").Value = '  -2.68%  '
$ws.Range("B49").Value = 'BabyDogeCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D49").Value = '0.0₆0222'
$ws.Range("E49").Value = '  +7.41%  '
$ws.Range("B50").Value = 'VeChain'
$ws.Range("C50").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0211'
$ws.Range("E50").Value = '  -1.27%  '
$ws.Range("B51").Value = 'WhiteBITCoin'
$ws.Range("C51").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '10.99'
$ws.Range("E51").Value = '  -0.47%  '